# Inclusion of Common IC Type information in 2010 data
#
# The sheet is a frequency table of surfaceWaterBodyIntercalibrationTypeCode
# values (column A) against their counts (column B). A new code,
# "RW-R-M4", needs to be inserted right after the existing "RW-R-M1" row
# (row 25), pushing every following category down by one row, and all of
# the counts need to be refreshed to reflect the updated 2010 dataset
# (which now also includes one more "blank code" row at the very bottom).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new "RW-R-M4" category: insert a blank row right after
# row 24 (RW-R-M1), shifting rows 25-32 down to 26-33.
$ws.Rows("25:25").Insert()

# New category label for the freshly inserted row.
$ws.Range("A25").Value = "RW-R-M4"

# Refresh every count in column B to the updated 2010 totals. Row numbers
# below are the FINAL (post-insert) row positions.
$counts = @{
    2  = 836
    3  = 184
    4  = 14
    5  = 25
    6  = 133
    7  = 41
    8  = 99
    9  = 533
    10 = 46
    11 = 2399
    12 = 59
    13 = 104
    14 = 291
    15 = 117
    16 = 163
    17 = 139
    18 = 43
    19 = 206
    20 = 143
    21 = 1038
    22 = 479
    23 = 323
    24 = 3
    25 = 1
    26 = 32
    27 = 13
    28 = 4
    29 = 202
    30 = 95
    31 = 1429
    32 = 3717
}

foreach ($row in $counts.Keys) {
    $ws.Cells.Item($row, 2).Value = $counts[$row]
}

# The inserted row pushed the former last row (an entry with a blank
# surfaceWaterBodyIntercalibrationTypeCode - column A stays empty) down to
# row 33; previously it held 5567 - now the underlying data yields 190 for
# that same blank-code bucket. Row 32 ("inapplicable") already got its own
# updated count above via $counts[32].
$ws.Range("B33").Value = 190
